# Commit: [Kadastro App] Kayıt silindi: 11483416
# The record with Kayıt No "11483416" was deleted from the workbook.
# It appears as a whole-row entry on both the master "Kayitlar" sheet
# (row 1689) and the filtered "Merkez İlçe" sheet (row 1150). Deleting
# the row shifts every following row up by one on each of those sheets.

$wb = $excel.ActiveWorkbook

$recordId = "11483416"

$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row
    $targetRow = -1

    for ($r = 2; $r -le $lastRow; $r++) {
        $val = $ws.Cells.Item($r, 1).Value2
        if ("$val" -eq $recordId) {
            $targetRow = $r
            break
        }
    }

    if ($targetRow -gt 0) {
        $ws.Rows.Item($targetRow).Delete()
    }
}
